{"js": "// Remove the \"Classification: Controlled\" text-box shapes that Word's\n// sensitivity-label add-in stamps into the first-page and even-page\n// footers. These are anchored (floating) shapes, not inline text, so we\n// reach them via Section.getFooter(...).shapes rather than Body text.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const section = sections.items[i];\n\n  // Primary footer never held the classification shape, but check all\n  // three footer kinds so the edit is robust regardless of section setup.\n  for (const kind of [\"FirstPage\", \"EvenPages\", \"Primary\"]) {\n    const footer = section.getFooter(kind);\n    const shapes = footer.shapes;\n    shapes.load(\"items\");\n    await context.sync();\n\n    if (shapes.items.length === 0) {\n      continue;\n    }\n\n    for (let j = 0; j < shapes.items.length; j++) {\n      shapes.items[j].body.load(\"text\");\n    }\n    await context.sync();\n\n    for (let j = shapes.items.length - 1; j >= 0; j--) {\n      const shape = shapes.items[j];\n      const shapeText = (shape.body && shape.body.text) || \"\";\n      if (shapeText.indexOf(\"Classification\") !== -1) {\n        shape.delete();\n      }\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the \"Classification: Controlled\" text-box shapes that Word's\n# sensitivity-label add-in stamped into the footers (anchored/floating\n# shapes, not part of the normal footer text run).\n$d = $word.ActiveDocument\n\nforeach ($sec in $d.Sections) {\n    # wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3\n    for ($idx = 1; $idx -le $sec.Footers.Count; $idx++) {\n        $ftr = $sec.Footers.Item($idx)\n        for ($j = $ftr.Shapes.Count; $j -ge 1; $j--) {\n            $shp = $ftr.Shapes.Item($j)\n            $shapeText = \"\"\n            if ($shp.TextFrame.HasText) {\n                $shapeText = $shp.TextFrame.TextRange.Text\n            }\n            if ($shapeText -match \"Classification\") {\n                $shp.Delete()\n            }\n        }\n    }\n}\n"}
